$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-27 Sunday" "2024-10-28 Monday"

Replace-Text "75×28=" "87×25="
Replace-Text "53×16=" "48×68="
Replace-Text "41×71=" "79×60="
Replace-Text "93×87=" "85×21="
Replace-Text "55×57=" "80×34="

Replace-Text "91×68=" "29×35="
Replace-Text "20×83=" "88×43="
Replace-Text "57×96=" "94×37="
Replace-Text "97×16=" "76×92="
Replace-Text "35×86=" "19×83="

Replace-Text "89×89=" "74×73="
Replace-Text "89×38=" "34×82="
Replace-Text "99×76=" "28×67="
Replace-Text "35×55=" "39×90="
Replace-Text "15×94=" "86×56="

Replace-Text "91×99=" "70×14="
Replace-Text "87×98=" "34×32="
Replace-Text "24×23=" "84×53="
Replace-Text "49×66=" "59×24="
Replace-Text "59×60=" "15×35="

Replace-Text "49×52=" "36×43="
Replace-Text "27×24=" "21×75="
Replace-Text "67×34=" "29×11="
Replace-Text "85×35=" "44×91="
Replace-Text "75×24=" "27×65="
